# Update statistical description table (rows 2-19) with refreshed
# dataset stats ("Commit with more data!") - counts grew and derived
# mean/std/percentile values shifted accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3106
$ws.Range("C2").Value = 477.1687057308436
$ws.Range("D2").Value = 50.46616697433213
$ws.Range("G2").Value = 454
$ws.Range("H2").Value = 515.75
$ws.Range("B3").Value = 8402
$ws.Range("C3").Value = 57.34298500357058
$ws.Range("D3").Value = 6.177880403290533
$ws.Range("F3").Value = 53.23
$ws.Range("G3").Value = 57.51
$ws.Range("H3").Value = 60.88
$ws.Range("B4").Value = 8402
$ws.Range("C4").Value = 11.28194715543918
$ws.Range("D4").Value = 15.89387772366126
$ws.Range("F4").Value = 6.73
$ws.Range("G4").Value = 11.74
$ws.Range("H4").Value = 13.42
$ws.Range("B5").Value = 8402
$ws.Range("C5").Value = 323.2435800999762
$ws.Range("D5").Value = 1.929981559305362
$ws.Range("F5").Value = 321.5625
$ws.Range("G5").Value = 323.52
$ws.Range("H5").Value = 324.89
$ws.Range("B6").Value = 8402
$ws.Range("C6").Value = 26.0857843370626
$ws.Range("D6").Value = 1.647480850675091
$ws.Range("F6").Value = 25.47
$ws.Range("G6").Value = 26.34
$ws.Range("H6").Value = 27.22
$ws.Range("B7").Value = 8402
$ws.Range("C7").Value = -46.57569626279457
$ws.Range("D7").Value = 10.18534554147162
$ws.Range("B8").Value = 8401
$ws.Range("C8").Value = 10.2717771693846
$ws.Range("D8").Value = 1.772277871327665
$ws.Range("B9").Value = 8402
$ws.Range("C9").Value = 9.385503451559153
$ws.Range("D9").Value = 1.676218270986525
$ws.Range("B10").Value = 8402
$ws.Range("C10").Value = 867.8385860509402
$ws.Range("D10").Value = 0.463405568739291
$ws.Range("B11").Value = 8401
$ws.Range("C11").Value = 1656.466492084276
$ws.Range("D11").Value = 1071.803536369196
$ws.Range("F11").Value = 750
$ws.Range("G11").Value = 1487
$ws.Range("H11").Value = 2532
$ws.Range("I11").Value = 3966
$ws.Range("B12").Value = 8402
$ws.Range("C12").Value = 1768.274696500833
$ws.Range("D12").Value = 1151.91236066656
$ws.Range("F12").Value = 794.25
$ws.Range("G12").Value = 1588.5
$ws.Range("H12").Value = 2688.75
$ws.Range("I12").Value = 4298
$ws.Range("B13").Value = 8402
$ws.Range("C13").Value = 0.5251831164008569
$ws.Range("D13").Value = 0.5409337532368425
$ws.Range("B14").Value = 8402
$ws.Range("C14").Value = 0.04610330873601525
$ws.Range("D14").Value = 0.02171281956609874
$ws.Range("B15").Value = 8402
$ws.Range("B16").Value = 8402
$ws.Range("B17").Value = 8402
$ws.Range("D17").Value = 0.00000000000005829017772356581
$ws.Range("B18").Value = 8402
$ws.Range("B19").Value = 8402
$ws.Range("C19").Value = 62.97569626279456
$ws.Range("D19").Value = 10.18534554147152
